# Reset the two sample/demo rows in the bed-info sheet so the uploaded
# template only contains an empty ("空床") bed entry instead of a
# checked-in student record, and point the active selection at C3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: building 14, dorm "E1224", bed 1, gender 女, status 空床 ---
$ws.Range("A2").Clear()
$ws.Range("B2").Value = 14
$ws.Range("C2").Value = "E1224"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "女"
$ws.Range("F2").Value = "空床"
$ws.Range("G2").Clear()
$ws.Range("H2").Clear()

# --- Row 3: building 14, dorm 1224 (numeric), bed 2, gender 女, status 空床 ---
$ws.Range("A3").Clear()
$ws.Range("B3").Value = 14
$ws.Range("C3").Value = 1224
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = "女"
$ws.Range("F3").Value = "空床"

# Move/save the active selection to C3, matching the saved workbook state.
[void]$ws.Range("C3").Select()
